$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the style of the existing
# header row (copy G1's formatting onto H1, then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Era / save-flag data for rows 2-32
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
